$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1762
$ws.Range("F9").Value = 314
$ws.Range("F10").Value = 1702
$ws.Range("F11").Value = 341
$ws.Range("F12").Value = 1406
$ws.Range("F14").Value = 324
$ws.Range("F15").Value = 666
$ws.Range("F16").Value = 12681
$ws.Range("F17").Value = 12704
$ws.Range("F18").Value = 943
$ws.Range("F19").Value = 735
$ws.Range("F21").Value = 501
$ws.Range("F23").Value = 518
$ws.Range("F27").Value = 236
$ws.Range("F28").Value = 666

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 50

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 157

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 157
$ws.Range("F6").Value = 1762
$ws.Range("F13").Value = 50
$ws.Range("F14").Value = 314
$ws.Range("F15").Value = 1702
$ws.Range("F16").Value = 341
$ws.Range("F17").Value = 1406
$ws.Range("F19").Value = 324
$ws.Range("F21").Value = 667
$ws.Range("F22").Value = 12681
$ws.Range("F23").Value = 12704
$ws.Range("F24").Value = 943
$ws.Range("F25").Value = 735
$ws.Range("F27").Value = 501
$ws.Range("F29").Value = 518
$ws.Range("F37").Value = 236
$ws.Range("F38").Value = 666
